$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "GenerateContract" column (F) to "Surcharge" — the data
# underneath (the literal "one" values) stays as-is.
$ws.Range("F1").Value = "Surcharge"

# Column G ("Last Name" header / "Test" values) becomes the new
# "GenerateContract" column, taking the same "one" value F used to carry,
# and picking up F1's (unstyled) formatting.
$ws.Range("G1").Value = "GenerateContract"
$ws.Range("G1").ClearFormats()

$ws.Range("G2:G7").Value = "one"

# Match the new column widths (F/G end up the same width) and the
# worksheet's lingering selection left over from the edit.
$ws.Columns("F").ColumnWidth = 28.33
$ws.Columns("G").ColumnWidth = 28.33

$ws.Range("F1").Select() | Out-Null
